$wb = $excel.ActiveWorkbook

$wsPub = $wb.Worksheets.Item("publications")
$wsPub.Range("C12").Value = 'J Integr Bioinform. 2021 Oct 22;18(3):20210026. doi: 10.1515/jib-2021-0026. PMID: 34674411; PMCID: PMC8573232. https://doi.org/10.1515/jib-2021-0026'
$wsPub.Range("A15").Value = 'Streptozotocin-Induced Diabetes Causes Changes in Serotonin-Positive Neurons in the Small Intestine in Pig Model'
$wsPub.Range("C22").Value = 'Examples and Counterexamples,  https://doi.org/10.1016/j.exco.2022.100087'
$wsPub.Range("A24").Value = 'Periportal steatosis in mice affects distinct parameters of pericentral drug metabolism '
$wsPub.Range("A25").Value = 'Uncertainty with Varying Subsurface Permeabilities Reduced Using Coupled Random Field and Extended Theory of Porous Media Contaminant Transport Models'
$wsPub.Range("A26").Value = 'Simulation of Contaminant Transport through the Vadose Zone: A Continuum Mechanical Approach within the Framework of the Extended Theory of Porous Media (eTPM)'
$wsPub.Range("A28").Value = 'Critical Evaluation of Discarded Donor Livers  in the Eurotransplant Region: Potential  Implications for Machine Perfusion'
$wsPub.Range("C28").Value = 'e-ISSN 2329-0358 © Ann Transplant, 2023; 28: e938132 DOI: 10.12659/AOT.938132'

$wsPres = $wb.Worksheets.Item("presentations")
$wsPres.Range("A5").Value = 'Fast 3D Isotropic High-Resolution MRI of Mouse Brain Using a Variable Flip Angle RARE Sequence With T2 Compensation @9.4T'
$wsPres.Range("A16").Value = 'Bayesian estimation shows the merits of reproducible and reusable modeling in systems biology'
$wsPres.Range("A22").Value = 'Multiscale and Multiphase Modeling of Function-Perfusion Processes in the Liver on Organ, Lobule and Cell Scale'
$wsPres.Range("A28").Value = 'Simulation Supported Liver Assessment for Donor Organs (SimLivA) - Con- tinuum-Biomechanical Modeling for Staging of Ischemia Reperfusion Injury During Liver Transplantation'

$wsOther = $wb.Worksheets.Item("other")
$wsOther.Range("A1").Value = 'Category'
$wsOther.Range("B1").Value = 'Title'
$wsOther.Range("C1").Value = 'Authors'
$wsOther.Range("A2").Value = 'Editoring'
$wsOther.Range("B2").Value = 'Computational Modeling for Liver Surgery and Interventions'
$wsOther.Range("C2").Value = 'Bruno Christ, Uta Dahmen, Nicole Radde, Tim Ricken'
$wsOther.Range("A3").Value = 'Organization Minisymposium'
$wsOther.Range("B3").Value = 'Computational Continuum Biomechanics'
$wsOther.Range("C3").Value = 'Tim Ricken, Oliver Röhrle, Silvia Budday'
$wsOther.Range("A4").Value = 'Organization Minisymposium'
$wsOther.Range("B4").Value = 'CONTINUUM BIOMECHANICS OF ACTIVE SYSTEMS'
$wsOther.Range("C4").Value = 'Tim Ricken, Oliver Röhrle, Silvia Budday'
$wsOther.Range("A5").Value = 'Editoring'
$wsOther.Range("B5").Value = 'Frontiers Research Topic “Multiscale Modeling for the Liver”'
$wsOther.Range("C5").Value = 'Ho H, Rezaina V, Schwen LO'
$wsOther.Range("A6").Value = 'Organization Minisymposium'
$wsOther.Range("B6").Value = 'Computational biomechanics and biomedical engineering of active biological systems – from methods to clinical application'
$wsOther.Range("C6").Value = 'Christian Bleiler, Lena Lambers, Renate Sachse'
$wsOther.Range("A7").Value = 'Editoring'
$wsOther.Range("B7").Value = 'Progress in Liver Stem Cell Therapy '
$wsOther.Range("C7").Value = 'Christ B, Oertel M '
$wsOther.Range("A8").Value = 'Accepted project proposal'
$wsOther.Range("B8").Value = 'Michael Stifel grant '
$wsOther.Range("C8").Value = 'Lena Lambers, Uta Dahmen '
$wsOther.Range("A9").Value = 'Workshop '
$wsOther.Range("B9").Value = 'Workshop on Computational Models in Biology and Medicine'
$wsOther.Range("C9").Value = 'Nicole Radde, Sebastian Höpfl '
$wsOther.Range("A10").Value = 'Miscellaneous'
$wsOther.Range("B10").Value = 'Kongresssekretär der Jahrestagung der mittelseutschen Viszeralmedizin'
$wsOther.Range("C10").Value = 'Hans-Michael Tautenhahn'
$wsOther.Range("A11").Value = 'Miscellaneous'
$wsOther.Range("B11").Value = ' Vorstandsmitglied der Thüringischen Gesellschaft für Chirurgie'
$wsOther.Range("C11").Value = 'Hans-Michael Tautenhahn'
$wsOther.Range("A12").Value = 'Miscellaneous'
$wsOther.Range("B12").Value = 'New Project Assistant'
$wsOther.Range("C12").Value = 'Hans-Michael Tautenhahn'
$wsOther.Range("A13").Value = 'Workshop Leader '
$wsOther.Range("B13").Value = 'GmdS/IBS Arbeitsgruppe Mathematical Models in Medicine and Biology'
$wsOther.Range("C13").Value = 'Nicole Radde, Ingmar Glauche '
$wsOther.Range("A14").Value = 'Accepted project proposal'
$wsOther.Range("B14").Value = 'Methodology for the calibration and analysis of stochastic models for heterogeneous intracellular processes with applications in cancer development'
$wsOther.Range("C14").Value = 'Nicole Radde'
$wsOther.Range("A15").Value = 'Accepted project proposal'
$wsOther.Range("B15").Value = 'Data-enhanced prediction of organ-specific tumor growth in the liver - a hybrid knowledge and data-driven approach'
$wsOther.Range("C15").Value = 'Tim Ricken'
$wsOther.Range("A16").Value = 'Editoring'
$wsOther.Range("B16").Value = 'Editor-in-Chief der Zeitschrift "Zeitschrift für Medizinische Physik" (Journal of Medical Physics) (ISSN 0939-3889)'
$wsOther.Range("C16").Value = 'Jürgen Reichenbach'
$wsOther.Range("A17").Value = 'Internship'
$wsOther.Range("B17").Value = 'Computational Modeling of Drug Detoxification – A Systems Medicine Approach'
$wsOther.Range("C17").Value = 'Matthias König'
$wsOther.Range("A18").Value = 'Accepted project proposal'
$wsOther.Range("B18").Value = 'SIMulation supported LIVer Assessment for donor organs (SimLivA) - Continuum-biomechanical modeling for staging of ischemia reperfusion injury during liver transplantation '
$wsOther.Range("C18").Value = 'Professorin Dr. Uta Dahmen, Professor Dr.-Ing. Tim Ricken, Privatdozent Dr. Hans-Michael Tautenhahn, Dr. Matthias König'
$wsOther.Range("A19").Value = 'Accepted project proposal'
$wsOther.Range("B19").Value = '031L0304X - CompLS - Runde 5 - Verbundprojekt: ATLAS - Al and Simulation for Tumor Liver Assessment - Entwicklung eines Systems zur klinischen Entscheidungsunterstützung in der Diagnose und Behandlung von Lebertumoren auf Basis von künstlicher Intelligenz und Simulationen '
$wsOther.Range("C19").Value = 'Professor Dr.-Ing. Tim Ricken, Professor Dr. Steffen Staab, Privatdozent Dr. Hans-Michael Tautenhahn, Dr. Matthias König'
$wsOther.Range("A20").Value = 'Editoring'
$wsOther.Range("B20").Value = 'Elected PEtab editor from March 2023 until March 2026'
$wsOther.Range("C20").Value = 'König Matthias'
$wsOther.Range("A21").Value = 'Accepted project proposal'
$wsOther.Range("B21").Value = 'X-Research Group: Physiologically based digital twins for the treatment of hypertension with ACE inhibitors and diuretics'
$wsOther.Range("C21").Value = 'König Matthias'
$wsOther.Range("A22").Value = 'Editoring'
$wsOther.Range("B22").Value = 'M. König has been elected as PEtab editor (2023-2026)'
$wsOther.Range("C22").Value = 'König Matthias'
$wsOther.Range("A23").Value = 'Dataset'
$wsOther.Range("B23").Value = 'Datasets for "Automated Detection of Portal Fields and Central Veins in Whole-Slide Images of Liver Tissue"'
$wsOther.Range("C23").Value = 'Budelmann D, Laue H, Weis N, Dahmen U, D’Alessandro LA, Biermayer I, Klingmüller U, Ghallab A, Hassan R, Begher-Tibbe B, Hengstler JG, Schwen LO'
$wsOther.Range("A24").Value = 'Dataset'
$wsOther.Range("B24").Value = 'Dataset for "Segmentation of Lipid Droplets in Histological Images"'
$wsOther.Range("C24").Value = 'Budelmann D, Cao Q, Laue H, Albadry M, Dahmen U, Schwen LO'
$wsOther.Range("A25").Value = 'Hosting and organization of COMBINE 2024'
$wsOther.Range("B25").Value = 'Hosting and organiztion of the Computational Modeling in Biology” Network (COMBINE) meeting 2024 in Stuttgart'
$wsOther.Range("C25").Value = 'Radde Nicole, Waltemath Dagmar, Höpfl Sebastian'
$wsOther.Range("A26").Value = 'Organization Minisymposium'
$wsOther.Range("B26").Value = 'In-silico Models of Coupled Biological Systems'
$wsOther.Range("C26").Value = 'Radde Nicole'
$wsOther.Range("A27").Value = 'Invited Speaker'
$wsOther.Range("B27").Value = 'Mathematische Modellbildung des Ischämie-Reperfusionsschadens zur Entscheidungsunterstützung bei Lebertransplantationen'
$wsOther.Range("C27").Value = 'Mandl L., Gerhäusser S., Lambers L., König M., Tautenhahn H.-M., Dahmen U., Ricken T.'
